$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 6238.8184
$ws.Range("I86").Value = 5824.75
$ws.Range("J86").Value = 6475.4287
$ws.Range("K86").Value = 5824.75
$ws.Range("L86").Value = 6475.4287
$ws.Range("M86").Value = -4701.75
$ws.Range("N86").Value = -8721.4287

# Row 88
$ws.Range("H88").Value = 2549.5
$ws.Range("I88").Value = 2499
$ws.Range("J88").Value = 2600
$ws.Range("K88").Value = 2499
$ws.Range("L88").Value = 2600
$ws.Range("M88").Value = -2093
$ws.Range("N88").Value = -3412

# Row 89
$ws.Range("H89").Value = 6238.8184
$ws.Range("I89").Value = 5824.75
$ws.Range("J89").Value = 6475.4287
$ws.Range("K89").Value = 29123.75
$ws.Range("L89").Value = 32377.1435
$ws.Range("M89").Value = -23507.75
$ws.Range("N89").Value = -43609.14350000001

# Row 91
$ws.Range("H91").Value = 2549.5
$ws.Range("I91").Value = 2499
$ws.Range("J91").Value = 2600
$ws.Range("K91").Value = 2499
$ws.Range("L91").Value = 2600
$ws.Range("M91").Value = -1095
$ws.Range("N91").Value = -5408

# Row 137
$ws.Range("H137").Value = 855.2222
$ws.Range("I137").Value = 682.8333
$ws.Range("K137").Value = 2048.4999
$ws.Range("M137").Value = 501.5001000000002

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4637.4375
$ws.Range("I32").Value = 4295
$ws.Range("K32").Value = 4295
$ws.Range("M32").Value = -4008

# Row 61
$ws.Range("H61").Value = 1791.5333
$ws.Range("I61").Value = 1791.5333
$ws.Range("K61").Value = 1791.5333
$ws.Range("M61").Value = -1579.5333

# Row 88
$ws.Range("H88").Value = 2414.1428
$ws.Range("I88").Value = 2399.6667
$ws.Range("J88").Value = 2425
$ws.Range("K88").Value = 2399.6667
$ws.Range("L88").Value = 2425
$ws.Range("M88").Value = -1993.6667
$ws.Range("N88").Value = -3237

# Row 91
$ws.Range("H91").Value = 2414.1428
$ws.Range("I91").Value = 2399.6667
$ws.Range("J91").Value = 2425
$ws.Range("K91").Value = 2399.6667
$ws.Range("L91").Value = 2425
$ws.Range("M91").Value = -995.6667000000002
$ws.Range("N91").Value = -5233

# Row 136
$ws.Range("H136").Value = 1791.5333
$ws.Range("I136").Value = 1791.5333
$ws.Range("K136").Value = 5374.5999
$ws.Range("M136").Value = -2824.5999

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1936.3334
$ws.Range("J20").Value = 2404.5
$ws.Range("L20").Value = 2404.5
$ws.Range("N20").Value = -2898.5

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

# Row 86
$ws.Range("H86").Value = 3836.25
$ws.Range("I86").Value = 4531.6665
$ws.Range("J86").Value = 1750
$ws.Range("K86").Value = 4531.6665
$ws.Range("L86").Value = 1750
$ws.Range("M86").Value = -3408.6665
$ws.Range("N86").Value = -3996

# Row 89
$ws.Range("H89").Value = 3836.25
$ws.Range("I89").Value = 4531.6665
$ws.Range("J89").Value = 1750
$ws.Range("K89").Value = 22658.3325
$ws.Range("L89").Value = 8750
$ws.Range("M89").Value = -17042.3325
$ws.Range("N89").Value = -19982

# Row 94
$ws.Range("H94").Value = 3023.8
$ws.Range("I94").Value = 2190.8462
$ws.Range("K94").Value = 2190.8462
$ws.Range("M94").Value = -1739.8462

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 344.42856
$ws.Range("I2").Value = 62.666668
$ws.Range("J2").Value = 555.75
$ws.Range("K2").Value = 62.666668
$ws.Range("L2").Value = 555.75
$ws.Range("M2").Value = 50.333332
$ws.Range("N2").Value = -781.75

# Row 132
$ws.Range("H132").Value = 4999.4443
$ws.Range("I132").Value = 4999.4443
$ws.Range("K132").Value = 14998.3329
$ws.Range("M132").Value = -12468.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 6378.4546
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 6966.3
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 6966.3
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -7556.3

# Row 27
$ws.Range("H27").Value = 6378.4546
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 6966.3
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 6966.3
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -7180.3

# Row 46
$ws.Range("H46").Value = 4636.273
$ws.Range("I46").Value = 3142.7144
$ws.Range("J46").Value = 5333.2666
$ws.Range("K46").Value = 3142.7144
$ws.Range("L46").Value = 5333.2666
$ws.Range("M46").Value = -2954.7144
$ws.Range("N46").Value = -5709.2666

# Row 68
$ws.Range("H68").Value = 4999.5
$ws.Range("I68").Value = 4999
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 4999
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -4250
$ws.Range("N68").Value = -6498

# Row 71
$ws.Range("H71").Value = 4999.5
$ws.Range("I71").Value = 4999
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 24995
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -21251
$ws.Range("N71").Value = -32488

# Row 82
$ws.Range("H82").Value = 774.75
$ws.Range("I82").Value = 750
$ws.Range("J82").Value = 799.5
$ws.Range("K82").Value = 750
$ws.Range("L82").Value = 799.5
$ws.Range("M82").Value = -389
$ws.Range("N82").Value = -1521.5

# Row 85
$ws.Range("H85").Value = 774.75
$ws.Range("I85").Value = 750
$ws.Range("J85").Value = 799.5
$ws.Range("K85").Value = 750
$ws.Range("L85").Value = 799.5
$ws.Range("M85").Value = 498
$ws.Range("N85").Value = -3295.5

# Row 93
$ws.Range("H93").Value = 663.0833
$ws.Range("I93").Value = 645.8
$ws.Range("K93").Value = 645.8
$ws.Range("M93").Value = 602.2

# Row 100
$ws.Range("H100").Value = 2799.6667
$ws.Range("I100").Value = 2499
$ws.Range("K100").Value = 2499
$ws.Range("M100").Value = -1958

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 8000
$ws.Range("J62").Value = 8000
$ws.Range("L62").Value = 8000
$ws.Range("N62").Value = -9248

# Row 65
$ws.Range("H65").Value = 8000
$ws.Range("J65").Value = 8000
$ws.Range("L65").Value = 40000
$ws.Range("N65").Value = -46240
